$d = $word.ActiveDocument
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>,.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>What</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t xml:space="preserve"> are the factors that affect IT Implementation?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:br w:type="page"/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>Introduction</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>As a practitioner, you are required to demonstrate that you have a clear understanding of the task you are performing. The purpose of the introduction is for the practitioner to demonstrate a clear understanding of the task at hand.</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>The introduction should include:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>What is the task?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>What is the purpose?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>What are you delivering, and what is the use of it?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>How is your report structured?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>Any other information you see as being important to include.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rStyle w:val="Heading1Char"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Heading1Char"/>
        </w:rPr>
        <w:t>Structure</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>A professional report at workplaces must consider a structure that helps readers understand discussions more effectively. Structure helps the flow of the contents and a logical connection between sections.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Look at your analysis outcome and the groups of themes or factors. It is a good place to think about shaping the main headings of the write-up.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rStyle w:val="Heading1Char"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Heading1Char"/>
        </w:rPr>
        <w:t>Methodology</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">As an IT professional, you must demonstrate how you get data and how you analyse it </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>t</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> support your discussions and recommendations.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">What are your sources of </w:t>
      </w:r>
      <w:r>
        <w:t>data?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>How did you find them? Where from?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">Why did you select these </w:t>
      </w:r>
      <w:r>
        <w:t>sources</w:t>
      </w:r>
      <w:r>
        <w:t>?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>What was your approach to the analysis of data?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>Any other important issue related to the methodology you would like to present.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/>
          <w:b/>
          <w:bCs/>
          <w:caps/>
          <w:spacing w:val="4"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:br w:type="page"/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>Outcome</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Explaining the outcome of the analysis</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Implement the structure and describe your analysis. Provide examples/evidence to justify that your work is accurate and based on data.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Delivery Factors</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>We have 4 different delivery factors to consider</w:t>
      </w:r>
      <w:r>
        <w:t>, turns out that these factors align perfectly with testing plan types.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Unit Test Plan</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Integration Test Plan</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>System Test Plan</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Acceptance Test Plan</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="1922983833"/>
          <w:citation/>
        </w:sdtPr>
        <w:sdtEndPr/>
        <w:sdtContent>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve">CITATION Sof \l 5129 </w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:noProof/>
            </w:rPr>
            <w:t>(Software Testing Fundamentals, n.d.)</w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:sdtContent>
      </w:sdt>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>These plans if executed correctly help prevent the fatal errors identified in the analysis.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Example</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>“</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>The objectives for the use of the Beta environment as a substitute for the pilot did not cover the primary in situ testing function that a pilot would typically perform</w:t>
      </w:r>
      <w:r>
        <w:t>”</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Ref: 3</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>A simple functional testing coverage would solve this inside a unit test.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Implementation Factors</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">A Project plan or </w:t>
      </w:r>
      <w:r>
        <w:t>a project schedule</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> would solve most if not all the issues. Each factor will be addressed with a different part of the plan.</w:t>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="-2060311414"/>
          <w:citation/>
        </w:sdtPr>
        <w:sdtEndPr/>
        <w:sdtContent>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve"> CITATION Sou \l 5129 </w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:noProof/>
            </w:rPr>
            <w:t xml:space="preserve"> (South Aftrican Goverment)</w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="1754015418"/>
          <w:citation/>
        </w:sdtPr>
        <w:sdtEndPr/>
        <w:sdtContent>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve"> CITATION MIC19 \l 5129 </w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:noProof/>
            </w:rPr>
            <w:t xml:space="preserve"> (ROBERTS, 2019)</w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:sdtContent>
      </w:sdt>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading3"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t xml:space="preserve">Factor: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Project Monitoring</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading4"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Administrative tasks</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Allocate time within the project plan to accommodate administrative tasks</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> this could include</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> status reports, team meetings</w:t>
      </w:r>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> etc.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>There needs to be regular updates or intervals to update the management on the project progression and feedback on that progress.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Example</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>“</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Found that Ministers were not always well served. Reporting to Ministers has been inconsistent, at times unduly optimistic and sometimes misrepresented the situation</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t xml:space="preserve">” Ref </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>29</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading3"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t xml:space="preserve">Factor: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Schedules/Deadlines</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading4"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Critical Path</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Critical path analysis to identify those tasks which are critical to the success and timely completion of the project.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> This then can be structed into something like milestones or deadlines.</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading3"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Factor: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Scope Creep</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading4"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Adjustment Plan</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading3"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t xml:space="preserve">Factor: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Technology Tools</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading4"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Project Approach</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading3"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t xml:space="preserve">Factor: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Workflow Communication</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading4"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
        <w:t>Communication Plan</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:lang w:eastAsia="en-NZ"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:br w:type="page"/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>Conclusion, summary, reflection</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Demonstrate your in-depth knowledge of the work you have completed and your confidence in pointing out the key points in a concise fashion.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Present a summary of what you did. Explain limitations. Outline what worked and what did not. Suggest how the outcome you delivered may be improved. Explain what your contribution to the IT community is? What is the value of your work?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
    </w:p>
  
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
$d.Content.InsertXML($xml)
Write-Output "done"
